$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line to drop the "EQD (10), " part.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EQD (10), EQN (12)", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EQN (12)", 2)

# 2. Append a new "Requisitos" section (Heading2) followed by a bullet
#    list item, after the final "Bibliografia" paragraph.
$last = $d.Paragraphs.Last
$end = $last.Range.End

$headingRange = $d.Range($end, $end)
$headingRange.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Style = "Heading2"
$headingPara.Range.Text = "Requisitos"

$bulletEnd = $d.Paragraphs.Last.Range.End
$bulletRange = $d.Range($bulletEnd, $bulletEnd)
$bulletRange.InsertParagraphAfter()
$bulletPara = $d.Paragraphs.Last
$bulletPara.Style = "ListBullet"
$bulletPara.Range.Text = "LOQ4044 -  Introdução à Engenharia da Qualidade  (Requisito fraco)`v"
